# Apply updated cryptocurrency price/volume data to Sheet1
# (values refreshed by the scheduled GitHub Actions scrape run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.072.00"
$ws.Range("E2").Value = "'  -1.20%  "
$ws.Range("D3").Value = "'1.793.18"
$ws.Range("E3").Value = "'  -0.06%  "
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("D5").Value = "'225.36"
$ws.Range("E5").Value = "'  +0.90%  "
$ws.Range("D6").Value = "'0.549"
$ws.Range("E6").Value = "'  -0.04%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "'  -0.08%  "
$ws.Range("D8").Value = "'32.52"
$ws.Range("E8").Value = "'  +1.31%  "
$ws.Range("E9").Value = "'  -0.96%  "
$ws.Range("D10").Value = "'0.0706"
$ws.Range("E10").Value = "'  +0.66%  "
$ws.Range("D11").Value = "'0.0930"
$ws.Range("E11").Value = "'  -0.05%  "
$ws.Range("D12").Value = "'2.050.86"
$ws.Range("E12").Value = "'  -0.14%  "
$ws.Range("D13").Value = "'1.809.91"
$ws.Range("E13").Value = "'  +1.03%  "
$ws.Range("E14").Value = "'  -1.13%  "
$ws.Range("E15").Value = "'  -1.97%  "
$ws.Range("D16").Value = "'34.033.88"
$ws.Range("E16").Value = "'  -1.33%  "
$ws.Range("D17").Value = "'4.17"
$ws.Range("E17").Value = "'  -2.07%  "
$ws.Range("D18").Value = "'68.03"
$ws.Range("E18").Value = "'  -1.05%  "
$ws.Range("D19").Value = "'243.31"
$ws.Range("E19").Value = "'  -2.40%  "
$ws.Range("E20").Value = "'  -1.15%  "
$ws.Range("E21").Value = "'  -0.07%  "
$ws.Range("D22").Value = "'10.68"
$ws.Range("E22").Value = "'  -2.76%  "
$ws.Range("D23").Value = "'4.09"
$ws.Range("E23").Value = "'  -2.80%  "
$ws.Range("E24").Value = "'  -2.04%  "
$ws.Range("D25").Value = "'159.11"
$ws.Range("E25").Value = "'  -1.31%  "
$ws.Range("D26").Value = "'16.26"
$ws.Range("E26").Value = "'  -0.22%  "
$ws.Range("D27").Value = "'7.01"
$ws.Range("E27").Value = "'  -0.91%  "
$ws.Range("E28").Value = "'  -1.21%  "
$ws.Range("E29").Value = "'  -0.04%  "
$ws.Range("D30").Value = "'0.0519"
$ws.Range("E30").Value = "'  -0.57%  "
$ws.Range("E31").Value = "'  +2.41%  "
$ws.Range("E32").Value = "'  -2.49%  "
$ws.Range("E33").Value = "'  -2.17%  "
$ws.Range("E34").Value = "'  -2.90%  "
$ws.Range("D35").Value = "'1.390.42"
$ws.Range("E35").Value = "'  -2.24%  "
$ws.Range("E36").Value = "'  +1.82%  "
$ws.Range("E37").Value = "'  -1.85%  "
$ws.Range("D38").Value = "'0.0185"
$ws.Range("E38").Value = "'  -2.41%  "
$ws.Range("E39").Value = "'  +0.03%  "
$ws.Range("D40").Value = "'79.15"
$ws.Range("E40").Value = "'  -5.24%  "
$ws.Range("D43").Value = "'2.18"
$ws.Range("E43").Value = "'  +1.91%  "
$ws.Range("D44").Value = "'0.0₆0142"
$ws.Range("E44").Value = "'  +14.14%  "
$ws.Range("E46").Value = "'  -0.60%  "
$ws.Range("D48").Value = "'5.87"
$ws.Range("E48").Value = "'  -2.03%  "
$ws.Range("D49").Value = "'1.950.07"
$ws.Range("E49").Value = "'  +0.21%  "

# Rows 41/42, 45/47, 50/51 swapped coin ordering; update all four columns.
$ws.Range("B41").Value = "'ARBITRUM"
$ws.Range("C41").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'0.915"
$ws.Range("E41").Value = "'  -3.31%  "
$ws.Range("B42").Value = "'MXToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.70"
$ws.Range("E42").Value = "'  -3.03%  "
$ws.Range("B45").Value = "'Quant"
$ws.Range("C45").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'107.70"
$ws.Range("E45").Value = "'  +2.15%  "
$ws.Range("B47").Value = "'Kaspa"
$ws.Range("C47").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").Value = "'0.0495"
$ws.Range("E47").Value = "'  -0.58%  "
$ws.Range("B50").Value = "'InjectiveProtocol"
$ws.Range("C50").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'12.00"
$ws.Range("E50").Value = "'  -1.17%  "
$ws.Range("B51").Value = "'PaxDollar"
$ws.Range("C51").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "'  -0.17%  "
